# Apply updated crypto price/volume data per diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some "D" price values look like plain decimals (e.g. "1.005") which Excel
# would otherwise auto-convert to a number on assignment. Force those cells to
# Text format first so the literal string is preserved, matching the source data
# (inline/shared string cells), then clear the temporary formatting afterward so
# the cell keeps its original (default) style.
$textForceCells = @(
    "D4",
    "D5",
    "D6",
    "D8",
    "D9",
    "D10",
    "D11",
    "D13",
    "D14",
    "D15",
    "D18",
    "D19",
    "D20",
    "D21",
    "D22",
    "D24",
    "D25",
    "D26",
    "D27",
    "D28",
    "D29",
    "D30",
    "D31",
    "D33",
    "D34",
    "D35",
    "D36",
    "D37",
    "D39",
    "D40",
    "D41",
    "D42",
    "D43",
    "D46",
    "D47",
    "D48",
    "D49",
    "D50",
    "D51"
)
foreach ($addr in $textForceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2
$ws.Range("D2").Value = "26.041.54"
$ws.Range("E2").Value = "  -0.11%  "

# Row 3
$ws.Range("D3").Value = "1.655.68"
$ws.Range("E3").Value = "  -0.56%  "

# Row 4
$ws.Range("D4").Value = "1.005"
$ws.Range("E4").Value = "  +0.07%  "

# Row 5
$ws.Range("D5").Value = "206.93"
$ws.Range("E5").Value = "  -1.05%  "

# Row 6
$ws.Range("D6").Value = "0.5170"
$ws.Range("E6").Value = "  +0.11%  "

# Row 7
$ws.Range("E7").Value = "  +0.03%  "

# Row 8
$ws.Range("D8").Value = "0.2575"
$ws.Range("E8").Value = "  -2.33%  "

# Row 9
$ws.Range("D9").Value = "0.06280"
$ws.Range("E9").Value = "  +1.12%  "

# Row 10
$ws.Range("D10").Value = "20.75"
$ws.Range("E10").Value = "  -0.63%  "

# Row 11
$ws.Range("D11").Value = "0.07560"
$ws.Range("E11").Value = "  +0.79%  "

# Row 12
$ws.Range("D12").Value = "1.677.14"
$ws.Range("E12").Value = "  +0.55%  "

# Row 13
$ws.Range("D13").Value = "4.385"

# Row 14
$ws.Range("D14").Value = "0.5373"
$ws.Range("E14").Value = "  -3.59%  "

# Row 15
$ws.Range("D15").Value = "66.19"
$ws.Range("E15").Value = "  +1.35%  "

# Row 16
$ws.Range("D16").Value = "0.0₅7914"
$ws.Range("E16").Value = "  -0.34%  "

# Row 17
$ws.Range("D17").Value = "26.056.94"
$ws.Range("E17").Value = "  -0.16%  "

# Row 18
$ws.Range("D18").Value = "1.005"
$ws.Range("E18").Value = "  +0.04%  "

# Row 19
$ws.Range("D19").Value = "4.685"
$ws.Range("E19").Value = "  -1.97%  "

# Row 20
$ws.Range("D20").Value = "187.30"
$ws.Range("E20").Value = "  +0.93%  "

# Row 21
$ws.Range("D21").Value = "10.05"
$ws.Range("E21").Value = "  -3.04%  "

# Row 22
$ws.Range("D22").Value = "6.152"
$ws.Range("E22").Value = "  +0.19%  "

# Row 23
$ws.Range("E23").Value = "  +0.06%  "

# Row 24
$ws.Range("D24").Value = "147.94"
$ws.Range("E24").Value = "  +1.24%  "

# Row 25
$ws.Range("D25").Value = "0.1211"
$ws.Range("E25").Value = "  -2.58%  "

# Row 26
$ws.Range("D26").Value = "7.345"
$ws.Range("E26").Value = "  -2.65%  "

# Row 27
$ws.Range("D27").Value = "15.64"
$ws.Range("E27").Value = "  -0.22%  "

# Row 28
$ws.Range("D28").Value = "1.396"
$ws.Range("E28").Value = "  +4.24%  "

# Row 29
$ws.Range("D29").Value = "0.05978"
$ws.Range("E29").Value = "  -5.26%  "

# Row 30
$ws.Range("D30").Value = "1.258"
$ws.Range("E30").Value = "  -1.00%  "

# Row 31
$ws.Range("D31").Value = "3.460"
$ws.Range("E31").Value = "  -0.23%  "

# Row 32
$ws.Range("E32").Value = "  -0.98%  "

# Row 33
$ws.Range("D33").Value = "1.630"
$ws.Range("E33").Value = "  +0.96%  "

# Row 34
$ws.Range("D34").Value = "0.9806"
$ws.Range("E34").Value = "  -1.23%  "

# Row 35
$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D35").Value = "2.390"
$ws.Range("E35").Value = "  -0.77%  "

# Row 36
$ws.Range("B36").Value = "MXToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D36").Value = "2.750"
$ws.Range("E36").Value = "  +1.82%  "

# Row 37
$ws.Range("D37").Value = "0.5855"
$ws.Range("E37").Value = "  -2.87%  "

# Row 38
$ws.Range("D38").Value = "1.096.73"
$ws.Range("E38").Value = "  +1.89%  "

# Row 39
$ws.Range("D39").Value = "0.01591"
$ws.Range("E39").Value = "  -0.80%  "

# Row 40
$ws.Range("D40").Value = "5.895"
$ws.Range("E40").Value = "  -2.87%  "

# Row 41
$ws.Range("D41").Value = "0.8445"
$ws.Range("E41").Value = "  -1.69%  "

# Row 42
$ws.Range("D42").Value = "1.003"
$ws.Range("E42").Value = "  -0.05%  "

# Row 43
$ws.Range("D43").Value = "100.22"
$ws.Range("E43").Value = "  +0.99%  "

# Row 44
$ws.Range("D44").Value = "1.813.33"
$ws.Range("E44").Value = "  +0.01%  "

# Row 45
$ws.Range("D45").Value = "0.0₈107"
$ws.Range("E45").Value = "  -0.89%  "

# Row 46
$ws.Range("D46").Value = "54.90"
$ws.Range("E46").Value = "  -1.94%  "

# Row 47
$ws.Range("D47").Value = "0.9943"
$ws.Range("E47").Value = "  -1.12%  "

# Row 48
$ws.Range("D48").Value = "7.997"
$ws.Range("E48").Value = "  +1.23%  "

# Row 49
$ws.Range("D49").Value = "0.05225"
$ws.Range("E49").Value = "  -0.48%  "

# Row 50
$ws.Range("D50").Value = "0.4239"
$ws.Range("E50").Value = "  -0.50%  "

# Row 51
$ws.Range("D51").Value = "5.842"
$ws.Range("E51").Value = "  -0.75%  "

# Restore default formatting on the cells we temporarily forced to Text
foreach ($addr in $textForceCells) {
    $ws.Range($addr).ClearFormats()
}
